$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.04792813956737518
$ws.Range("C2").Value = 0.02637087181210518
$ws.Range("D2").Value = 0.014985552988946438
$ws.Range("E2").Value = 0.006959357298910618
$ws.Range("F2").Value = 0.0008878824883140624
$ws.Range("G2").Value = 0.0011897289659827948
$ws.Range("J2").Value = 0.12786182761192322
$ws.Range("K2").Value = 1.4881341457366943
